# openFAST_config.xlsx update:
#  - split template paths from generated output paths
#  - shrink example simulation duration for faster validation
#  - new IEC_template.IPT (S-class turbulence support)
#  - active sheet/selection moved to DLC_List

$wb = $excel.ActiveWorkbook

$config = $wb.Worksheets.Item("config")
$dlc    = $wb.Worksheets.Item("DLC_List")

# --- config sheet: template file locations now live under examples/templates ---
$config.Cells.Item(2, 2).Value2  = "../examples/templates/sim/5MW_Land_DLL_WTurb.fst"
$config.Cells.Item(3, 2).Value2  = "../examples/templates/sim/NRELOffshrBsline5MW_Onshore_ElastoDyn.dat"
$config.Cells.Item(4, 2).Value2  = "../examples/templates/sim/NRELOffshrBsline5MW_Onshore_ServoDyn.dat"
$config.Cells.Item(5, 2).Value2  = "../examples/templates/sim/NRELOffshrBsline5MW_Onshore_AeroDyn15.dat"
$config.Cells.Item(6, 2).Value2  = "../examples/templates/sim/NRELOffshrBsline5MW_InflowWind_12mps.dat"
$config.Cells.Item(7, 2).Value2  = "../examples/templates/wind/TurbSim.inp"
$config.Cells.Item(8, 2).Value2  = "../examples/templates/wind/IEC_template.IPT"

# --- config sheet: generated output now lives under examples/generated ---
$config.Cells.Item(10, 2).Value2 = "../examples/generated/sim"
$config.Cells.Item(11, 2).Value2 = "../examples/generated/wind"

# --- DLC_List sheet: reduce example simulation Duration 60 -> 10 seconds ---
$dlc.Cells.Item(2, 7).Value2 = "10"
$dlc.Cells.Item(3, 7).Value2 = "10"
$dlc.Cells.Item(4, 7).Value2 = "10"
$dlc.Cells.Item(5, 7).Value2 = "10"

# --- selection / active sheet bookkeeping ---
$config.Range("B8").Select() | Out-Null
$dlc.Range("G6").Select() | Out-Null
$dlc.Activate() | Out-Null
